# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest pull from coinranking.com (scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.645.27"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "2.655.55"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "597.76"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "2.655.73"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E13").Value = "  -1.89%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "28.09"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "3.137.87"
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("E16").Value = "  -3.38%  "
$ws.Range("D17").Value = "67.558.56"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "2.650.85"
$ws.Range("E18").Value = "  -2.58%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "8.33"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +7.93%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "12.06"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "363.99"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.64%  "
$ws.Range("E22").Value = "  -2.57%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.81"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.50%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.04"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +8.43%  "
$ws.Range("E25").Value = "  -4.37%  "
$ws.Range("E26").Value = "  +0.09%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "71.03"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("D28").Value = "2.795.17"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("E30").Value = "  -0.15%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "558.62"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -5.79%  "
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -5.07%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "158.41"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.75%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "19.45"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  -5.43%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "40.27"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("E48").Value = "  -4.19%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "154.84"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "3.90"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("E51").Value = "  -3.14%  "
